$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.714.95'
$ws.Range("E2").Value = '  +0.67%  '
$ws.Range("D3").Value = '3.625.85'
$ws.Range("E3").Value = '  +1.51%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '611.32'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.58'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.58%  '
$ws.Range("D7").Value = '3.626.37'
$ws.Range("E7").Value = '  +1.56%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("E9").Value = '  -0.61%  '
$ws.Range("E10").Value = '  +0.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.95'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.32%  '
$ws.Range("E12").Value = '  +0.76%  '
$ws.Range("D13").Value = '4.245.47'
$ws.Range("E13").Value = '  +1.71%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '30.08'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.20%  '
$ws.Range("D16").Value = '3.625.25'
$ws.Range("E16").Value = '  +1.28%  '
$ws.Range("D17").Value = '66.855.11'
$ws.Range("E17").Value = '  +0.76%  '
$ws.Range("E18").Value = '  +1.62%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.65'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.32%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.39'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.13'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.57%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '429.06'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.26%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.621'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.92%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.85'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.45%  '
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000124'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.58%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.45'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.61'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.53'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.64%  '
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("D31").Value = '3.628.46'
$ws.Range("E31").Value = '  +1.85%  '
$ws.Range("E32").Value = '  +1.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.160'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.11%  '
$ws.Range("E34").Value = '  -0.58%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.92'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.48%  '
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("E37").Value = '  +1.48%  '
$ws.Range("E38").Value = '  -1.62%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '176.78'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.55%  '
$ws.Range("E40").Value = '  +1.74%  '
$ws.Range("E41").Value = '  +0.63%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.902'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.41%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.90'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.98%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '46.26'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.62'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +8.91%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.13%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '25.23'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.23%  '
$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.18'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.42%  '
$ws.Range("E49").Value = '  +1.88%  '
$ws.Range("E50").Value = '  +1.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.965'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.92%  '
